$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1) Remove the horizontal-rule (w:pict) content from its paragraph, leaving it empty.
Set-ParagraphXml 4 '<w:p w14:paraId="40401C1E" w14:textId="77777777" w:rsidR="00437283" w:rsidRPr="00437283" w:rsidRDefault="00437283" w:rsidP="00437283"/>'

# 2) Split "framework" into its own run, wrapped in proofErr spell-check markers.
Set-ParagraphXml 5 ('<w:p w14:paraId="76838615" w14:textId="77777777" w:rsidR="00437283" w:rsidRDefault="00437283" w:rsidP="00437283">' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve">Para la solución de este desafío, se plantearon varias alternativas, las cuales fueron posteriormente implementadas en el entorno de trabajo </w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Qt</w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve">. Este </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r w:rsidRPr="00437283"><w:t>framework</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve"> nos permitió, mediante el uso responsable de sus librerías, manipular imágenes de manera eficiente, además de exportarlas y guardarlas en archivos</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
    '</w:p>')

# 3) Split "char" into its own run, wrapped in proofErr spell-check markers.
Set-ParagraphXml 11 ('<w:p w14:paraId="770DC04C" w14:textId="77777777" w:rsidR="00437283" w:rsidRPr="00437283" w:rsidRDefault="00437283" w:rsidP="00437283">' + `
    '<w:r w:rsidRPr="00437283"><w:t>La información de las imágenes fue almacenada en un arreglo de caracteres (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r w:rsidRPr="00437283"><w:t>char</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve">), aprovechando que cada carácter ocupa 8 bits, lo que facilita el tratamiento a nivel de bits. Gracias a las operaciones explicadas en clase, desarrollamos funciones específicas para trabajar directamente sobre los bits de cada componente de color. Estas funciones incluyen operaciones como </w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>XOR</w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>rotaciones</w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:t xml:space="preserve"> y </w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>sumas de píxeles</w:t></w:r>' + `
    '<w:r w:rsidRPr="00437283"><w:t>, las cuales permiten desenmascarar las imágenes, liberándolas de las máscaras aplicadas.</w:t></w:r>' + `
    '</w:p>')

# 4) Add <w:semiHidden/> to the "Fuentedeprrafopredeter" (Default Paragraph Font) style.
$styles = $d.Styles
$style = $styles.Item("Fuentedeprrafopredeter")
$style.Font.Hidden = $true

Write-Host "done"
